$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new columns (C and D) before the existing "% Complete Actual"
#    column. This shifts the old column C (percent-complete-actual data) to
#    column E, and - importantly - this engine expands the existing A1:C1
#    merge to A1:E1 automatically and copies styles from the left neighbour
#    without creating duplicate style records.
# ---------------------------------------------------------------------------
$ws.Range("C1:D1").EntireColumn.Insert(-4161)

# ---------------------------------------------------------------------------
# 2. Title text (row 1)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "AWS Certified Developer Course"

# ---------------------------------------------------------------------------
# 3. Header row (row 3)
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Goal % Complete"
$ws.Range("B3").Value = "Goal Date"
$ws.Range("C3").Value = "Goal Section"
$ws.Range("D3").Value = "Actual Section Complete"
$ws.Range("E3").Value = "Actual % Complete"

# The newly inserted columns C/D copied their formatting from column B (the
# date column), so column C currently looks like a date column. Re-base C's
# formatting on column A's plain (non-date) header/body style instead.
$ws.Range("A3").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("A4:A14").Copy()
$ws.Range("C4:C14").PasteSpecial(-4122)
$ws.Range("D4:D14").PasteSpecial(-4122)

# D3 (Actual Section Complete header) and D4:D14 (body) need an integer
# ("0") number format while keeping the already-copied fill/border.
$ws.Range("D3").NumberFormat = "0"
$ws.Range("D4:D14").NumberFormat = "0"

# ---------------------------------------------------------------------------
# 4. Body data for the new "Goal Section" (C) and "Actual Section Complete"
#    (D, first row only) columns.
# ---------------------------------------------------------------------------
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 3
$ws.Range("C6").Value = 5
$ws.Range("C7").Value = 7
$ws.Range("C8").Value = 9.84
$ws.Range("C9").Value = 11.106
$ws.Range("C10").Value = 12
$ws.Range("C11").Value = 14.148
$ws.Range("C12").Value = 16
$ws.Range("C13").Value = 18
$ws.Range("C14").Value = 22

$ws.Range("D4").Value = 2

# ---------------------------------------------------------------------------
# 5. Course info bullet list (rows 17-24) and the extra note in E17.
# ---------------------------------------------------------------------------
$ws.Range("A17").Value = "Skill level: All Levels"
$ws.Range("A18").Value = "Students: 82452"
$ws.Range("A19").Value = "Languages: English"
$ws.Range("A20").Value = "Captions: Yes"
$ws.Range("A21").Value = "Practice tests: 1"
$ws.Range("A22").Value = "Questions: 65"
$ws.Range("A23").Value = "Lectures: 221"
$ws.Range("A24").Value = "Video: 19 total hours"

$ws.Range("E17").Value = " - need approx 2hrs of video per week to obtain goal"

# ---------------------------------------------------------------------------
# 6. Column widths for the new columns.
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 22.140625
$ws.Columns("D").ColumnWidth = 24.5703125

# ---------------------------------------------------------------------------
# 7. Selection cosmetic state.
# ---------------------------------------------------------------------------
$ws.Range("C16").Select()

Write-Host "stage1 done"
